# Fruta / hortaliza, semanal
# Insert a new weekly data row into the "Camote" block of the Zapallo
# (Vega Monumental Concepción) sheet. Inserting a whole row at position 87
# shifts all the existing rows 87-184 down to 88-185 (the row that used to
# be last, row 184, becomes row 185), and the sheet dimension grows from
# A1:R184 to A1:R185 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87..184 down to 88..185, leaving a blank row 87 to fill in.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly record.
$ws.Cells.Item(87, 1).Value = 11
$ws.Cells.Item(87, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(87, 3).Value = "Bíobío"
$ws.Cells.Item(87, 4).Value = 44629
$ws.Cells.Item(87, 5).Value = 8
$ws.Cells.Item(87, 6).Value = 100112045
$ws.Cells.Item(87, 7).Value = "Zapallo"
$ws.Cells.Item(87, 8).Value = "Camote"
$ws.Cells.Item(87, 9).Value = "1a (cosecha)"
$ws.Cells.Item(87, 10).Value = 500
$ws.Cells.Item(87, 11).Value = 300
$ws.Cells.Item(87, 12).Value = 350
$ws.Cells.Item(87, 13).Value = 330
$ws.Cells.Item(87, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(87, 15).Value = "Región del Maule"
$ws.Cells.Item(87, 16).Value = 330
$ws.Cells.Item(87, 17).Value = 1
$ws.Cells.Item(87, 18).Value = "Hortaliza"
